# Issue #5: property category labels were mis-tagged as "land" on the
# 建物 (building) and 汽車 (car) sheets. Correct the property_category
# column (row 2) on each sheet to match its actual sheet type.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: property_category column is I
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"

# 汽車 (car) sheet: property_category column is H
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
